$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C width (closest achievable to 9.125 given engine's width quantization) ---
$ws.Columns.Item(3).ColumnWidth = 8.14

# --- Fill in the two new cells on existing row 10 (PriceChange / UpDown) ---
$ws.Cells.Item(10, 24).Value = 0.6499990000000011
$ws.Cells.Item(10, 25).Value = "Up"

# --- Add new row 11 of trade data ---
$ws.Cells.Item(11, 1).Value = 42654.883275462962
$ws.Cells.Item(11, 2).Value = 27
$ws.Cells.Item(11, 3).Value = "Strong Buy"
$ws.Cells.Item(11, 4).Value = 40
$ws.Cells.Item(11, 5).Value = 8481
$ws.Cells.Item(11, 6).Value = 451
$ws.Cells.Item(11, 7).Value = 64
$ws.Cells.Item(11, 8).Value = 34
$ws.Cells.Item(11, 9).Value = 96
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(11, 11).Value = 16335
$ws.Cells.Item(11, 12).Value = 113
$ws.Cells.Item(11, 13).Value = 61
$ws.Cells.Item(11, 14).Value = 29
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(11, 16).Value = "Named"
$ws.Cells.Item(11, 17).Value = 58.438771163779279
$ws.Cells.Item(11, 18).Value = 0.49
$ws.Cells.Item(11, 19).Value = 0.093299999999999994
$ws.Cells.Item(11, 19).NumberFormat = "0.00%"
$ws.Cells.Item(11, 20).Value = 0.024899999999999999
$ws.Cells.Item(11, 20).NumberFormat = "0.00%"
$ws.Cells.Item(11, 21).Value = 2.34
$ws.Cells.Item(11, 22).Value = "N/A"
$ws.Cells.Item(11, 23).Value = 2
